$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new test row (row 13) -----------------------------------
# compare | qualifier | qualifer_term | identifier | term | pulmonary_2010 | pulmonary_2016 | transplant_2016
# (qualifer_term/identifier/term typed first, then compare/qualifier,
# then the three Y/N flag columns.)
$ws.Range("C13").Value = "disorders of environmental origin"
$ws.Range("D13").Value = "C21.866.915.300.200.150"
$ws.Range("E13").Value = "Brain Concussion"
$ws.Range("A13").Value = "Old only"
$ws.Range("B13").Value = "C21"
$ws.Range("F13").Value = "N"
$ws.Range("G13").Value = "N"
$ws.Range("H13").Value = "N"

# Give the new row the same border/format as the row above it ...
$ws.Range("A12:H12").Copy()
$ws.Range("A13:H13").PasteSpecial(-4122)   # xlPasteFormats
# ... but the two middle cells were typed in without that formatting,
# so they keep the plain/default (borderless) look.
$ws.Range("D13:E13").Borders.LineStyle = -4142   # xlLineStyleNone

# --- Tidy up the header row (row 2) so it stops using the red font ----
# Row 11 used a red warning font; clear it back to plain black text ...
$ws.Range("A11:H11").Font.ColorIndex = 1
# ... and bring the second header row (row 2) into line with it.
$ws.Range("A11:H11").Copy()
$ws.Range("A2:H2").PasteSpecial(-4122)   # xlPasteFormats

# --- Column D needs to be a bit wider for the long new identifier -----
$ws.Columns("D").ColumnWidth = 19.6

# --- Selection ends up on the newly entered cell -----------------------
[void]$ws.Range("B13").Select()
